# Update countries & provincias Spain
#
# Applies the COVID-19 data refresh captured in the commit:
#  - updates the "Datos actualizados..." timestamp
#  - refreshes total/new/active/recovered/critical/death counts for several
#    countries, which also re-sorts a handful of rows (by total cases,
#    column B, descending) causing country labels to shift between rows
#    62/63, 86/87/88 and 213/214.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 18 de Agosto de 2020 a las 05:22"

# --- Per-row updates ----------------------------------------------------
# Each entry maps a row number to the column letter / new value pairs that
# must be written. Country name (column A) is only included when the row's
# label changes because of the re-sort.
$updates = @{
    9   = @{ B = 541493; D = 370717; E = 144295; H = 26481 }                      # Peru
    29  = @{ B = 103300; C = 267;    D = 84445;  E = 17586 }                      # Kazajistan
    40  = @{ B = 78534;  C = 211;    D = 18003;  E = 50587; G = 5; H = 9944 }      # Belgica
    51  = @{ B = 50995;  C = 493;    D = 7450;   E = 41962; G = 8; H = 1583 }      # Honduras
    62  = @{ A = "Venezuela";      B = 34802; D = 23575; E = 10939; H = 288 }
    63  = @{ A = "Azerbaiyan";     B = 34343; D = 32042; E = 1793;  H = 508 }
    71  = @{ B = 23772;  C = 213;              E = 8795;  G = 17; H = 438 }       # Australia
    86  = @{ A = "Paraguay";       B = 10135; D = 6210;  E = 3780; H = 145 }
    87  = @{ A = "Noruega";        B = 10060; D = 8857;  E = 942;  H = 261 }
    88  = @{ A = "Zambia";         B = 9839;  D = 8575;  E = 1000; H = 264 }
    158 = @{ E = 491; G = 1; H = 25 }                                             # Vietnam
    213 = @{ A = "Islas Malvinas"; D = 13; H = 0 }
    214 = @{ A = "Montserrat";     D = 12; H = 1 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
